$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date column (week of Jan 15, next day 21_01_2024) with worked-day counts
$ws.Range("E1").Value = "21_01_2024"
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 4

# Move active selection to E6, matching the author's last selected cell
$ws.Range("E6").Select()
